$d = $word.ActiveDocument

# 1. Title paragraph: merge the two runs "Samenvatting stage " + "Ometa"
#    (which were split apart by a spell-check proofErr pair) into a single
#    run "Samenvatting stage Ometa", and drop the now-unneeded
#    <w:proofErr .../> spellStart/spellEnd markers entirely.
#
#    Range.InsertXML() on a plain-text range tends to leave a stray
#    <w:proofErr/> marker clinging to the boundary of the replaced range,
#    so instead we target a range that fully spans the title paragraph
#    *and* the (untouched) empty paragraph that follows it. Reproducing
#    both paragraphs verbatim (using the original paragraph/run
#    properties and ids) guarantees the proofErr markers -- which sit
#    strictly inside that span -- are dropped, while nothing else in the
#    document changes.
$titlePara = $d.Paragraphs(1)
$nextPara = $d.Paragraphs(2)
$titleRange = $d.Range($titlePara.Range.Start, $nextPara.Range.End)

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' +
            '<w:p w14:paraId="3B58F8B1" w14:textId="3910AD36" w:rsidR="00714D72" w:rsidRDefault="00570143" w:rsidP="00570143">' +
              '<w:pPr><w:pStyle w:val="Titel"/><w:jc w:val="center"/><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr>' +
              '<w:r w:rsidRPr="00052E81"><w:rPr><w:lang w:val="nl-BE"/></w:rPr><w:t>Samenvatting stage Ometa</w:t></w:r>' +
            '</w:p>' +
            '<w:p w14:paraId="2D346220" w14:textId="77777777" w:rsidR="004A6C07" w:rsidRPr="004A6C07" w:rsidRDefault="004A6C07" w:rsidP="004A6C07">' +
              '<w:pPr><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$titleRange.InsertXML($titleXml)

# 2. Introduce the typo: "document" -> "ocument" in the intro paragraph.
$d.Content.Find.Execute("Dit document vat", $true, $false, $false, $false, $false, $true, 1, $false, "Dit ocument vat", 2)
